# Updates the cryptocurrency price/volume table with refreshed figures.
# Cells whose new text looks like a plain decimal number (e.g. "303.41")
# are written with a leading apostrophe so Excel keeps them as text
# (matching the workbook's original inline-string/text cells) instead of
# silently converting them to floating point numbers; the cell style is
# then reset to "Normal" so no stray formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value2 = '41.834.10'
$ws.Range('E2').Value2 = '  -0.03%  '
$ws.Range('D3').Value2 = '2.271.56'
$ws.Range('E3').Value2 = '  +0.01%  '
$ws.Range('E4').Value2 = '  +0.04%  '
$ws.Range('D5').Value2 = '''303.41'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value2 = '  +0.06%  '
$ws.Range('D6').Value2 = '''92.82'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value2 = '  -0.25%  '
$ws.Range('E7').Value2 = '  +1.27%  '
$ws.Range('E8').Value2 = '  -0.01%  '
$ws.Range('D9').Value2 = '''0.485'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value2 = '  -0.70%  '
$ws.Range('D10').Value2 = '''32.54'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value2 = '  -0.03%  '
$ws.Range('D11').Value2 = '''53.37'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value2 = '  -2.07%  '
$ws.Range('E12').Value2 = '  -0.44%  '
$ws.Range('E13').Value2 = '  -1.59%  '
$ws.Range('D14').Value2 = '''6.69'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value2 = '  +0.34%  '
$ws.Range('D15').Value2 = '2.624.34'
$ws.Range('E15').Value2 = '  +0.15%  '
$ws.Range('D16').Value2 = '''14.30'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value2 = '  +0.59%  '
$ws.Range('D17').Value2 = '2.263.13'
$ws.Range('E17').Value2 = '  -1.88%  '
$ws.Range('D18').Value2 = '''0.782'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value2 = '  +3.57%  '
$ws.Range('D19').Value2 = '41.780.60'
$ws.Range('E19').Value2 = '  +0.14%  '
$ws.Range('D20').Value2 = '''12.81'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value2 = '  +3.66%  '
$ws.Range('E21').Value2 = '  -0.17%  '
$ws.Range('E22').Value2 = '  +0.05%  '
$ws.Range('E23').Value2 = '  +0.14%  '
$ws.Range('D24').Value2 = '''244.00'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value2 = '  +1.21%  '
$ws.Range('D25').Value2 = '''2.57'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value2 = '  -0.77%  '
$ws.Range('E26').Value2 = '  +2.84%  '
$ws.Range('E27').Value2 = '  -0.05%  '
$ws.Range('D28').Value2 = '''24.02'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value2 = '  +0.62%  '
$ws.Range('E29').Value2 = '  -1.54%  '
$ws.Range('D30').Value2 = '''2.07'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value2 = '  -5.72%  '
$ws.Range('D31').Value2 = '''34.96'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value2 = '  +2.12%  '
$ws.Range('D32').Value2 = '''160.30'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value2 = '  +1.28%  '
$ws.Range('E33').Value2 = '  +1.27%  '
$ws.Range('E34').Value2 = '  +0.01%  '
$ws.Range('D35').Value2 = '''0.0742'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value2 = '  +0.37%  '
$ws.Range('D36').Value2 = '''3.02'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value2 = '  -1.46%  '
$ws.Range('B37').Value2 = 'Celestia'
$ws.Range('C37').Value2 = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D37').Value2 = '''16.89'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value2 = '  +0.86%  '
$ws.Range('B38').Value2 = 'Kaspa'
$ws.Range('C38').Value2 = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').Value2 = '''0.106'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value2 = '  +1.80%  '
$ws.Range('E39').Value2 = '  +0.10%  '
$ws.Range('E40').Value2 = '  +0.61%  '
$ws.Range('E41').Value2 = '  +0.47%  '
$ws.Range('D42').Value2 = '''3.94'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value2 = '  -1.17%  '
$ws.Range('B43').Value2 = 'Maker'
$ws.Range('C43').Value2 = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value2 = '2.017.53'
$ws.Range('E43').Value2 = '  -2.15%  '
$ws.Range('B44').Value2 = 'EnergySwap'
$ws.Range('C44').Value2 = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').Value2 = '''19.51'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value2 = '  -3.29%  '
$ws.Range('E45').Value2 = '  +1.05%  '
$ws.Range('D46').Value2 = '''10.44'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value2 = '  +3.75%  '
$ws.Range('E47').Value2 = '  +6.78%  '
$ws.Range('E48').Value2 = '  -2.43%  '
$ws.Range('D49').Value2 = '''53.38'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value2 = '  +3.06%  '
$ws.Range('D50').Value2 = '''73.16'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value2 = '  +3.64%  '
$ws.Range('B51').Value2 = 'Stacks'
$ws.Range('C51').Value2 = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D51').Value2 = '''1.51'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value2 = '  -1.27%  '
